$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM results for the Tgfb3-Eng ligand/receptor sheet.
# Each array holds the new values for columns E..T (16 values) of the
# given row, in column order:
# E  Ligand-expressing cells
# F  Ligand detection rate
# G  Ligand average expression value
# H  Ligand total expression value
# I  Ligand derived specificity of average expression value
# J  Ligand derived specificity of total expression value
# K  Receptor-expressing cells
# L  Receptor detection rate
# M  Receptor average expression value
# N  Receptor total expression value
# O  Receptor derived specificity of average expression value
# P  Receptor derived specificity of total expression value
# Q  Edge average expression weight
# R  Edge total expression weight
# S  Edge average expression derived specificity
# T  Edge total expression derived specificity
$data = @{
    2 = @(3, 1, 2.29903, 6.89709, 0.04075801785348079, 0.04075801785348079, 3, 1, 201.098592, 603.295776, 0.7918622805845071, 0.791862280584507, 462.3316959657601, 4160.98526369184, 0.03227473696956135, 0.03227473696956135)
    3 = @(3, 1, 2.29903, 6.89709, 0.04075801785348079, 0.04075801785348079, 3, 1, 35.924535, 107.773605, 0.1414593902976603, 0.1414593902976603, 82.59158370105, 743.3242533094501, 0.005765604355294544, 0.005765604355294545)
    4 = @(3, 1, 2.29903, 6.89709, 0.04075801785348079, 0.04075801785348079, 3, 1, 16.93339666666667, 50.80019, 0.0666783291178327, 0.06667832911783268, 38.93038693856667, 350.3734824471, 0.002717676528624893, 0.002717676528624893)
    5 = @(3, 1, 22.27635266666667, 66.829058, 0.3949230674234065, 0.3949230674234066, 3, 1, 201.098592, 603.295776, 0.7918622805845071, 0.791862280584507, 4479.743156162112, 40317.68840545901, 0.3127246808253278, 0.3127246808253278)
    6 = @(3, 1, 22.27635266666667, 66.829058, 0.3949230674234065, 0.3949230674234066, 3, 1, 35.924535, 107.773605, 0.1414593902976603, 0.1414593902976603, 800.26761104601, 7202.40849941409, 0.05586557633219687, 0.05586557633219688)
    7 = @(3, 1, 22.27635266666667, 66.829058, 0.3949230674234065, 0.3949230674234066, 3, 1, 16.93339666666667, 50.80019, 0.0666783291178327, 0.06667832911783268, 377.2143159912245, 3394.92884392102, 0.02633281026588193, 0.02633281026588193)
    8 = @(3, 1, 31.831433, 95.494299, 0.5643189147231126, 0.5643189147231126, 3, 1, 201.098592, 603.295776, 0.7918622805845071, 0.791862280584507, 6401.256357642337, 57611.30721878102, 0.446862862789618, 0.4468628627896178)
    9 = @(3, 1, 31.831433, 95.494299, 0.5643189147231126, 0.5643189147231126, 3, 1, 35.924535, 107.773605, 0.1414593902976603, 0.1414593902976603, 1143.529428908655, 10291.76486017789, 0.07982820961016884, 0.07982820961016884)
    10 = @(3, 1, 31.831433, 95.494299, 0.5643189147231126, 0.5643189147231126, 3, 1, 16.93339666666667, 50.80019, 0.0666783291178327, 0.06667832911783268, 539.0142814574233, 4851.12853311681, 0.03762784232332586, 0.03762784232332585)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, 5 + $i).Value = $vals[$i]
    }
}

Write-Host "updated $($data.Keys.Count) rows"
